$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D16").Value = "2016-02-29 04:26:36"
$wsZh.Range("G16").Value = "2016-02-29 04:27:23"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D16").Value = "2016-02-29 04:26:47"
$wsDe.Range("G16").Value = "2016-02-29 04:27:43"
